$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / URL / Notes updates ---
$ws.Range("D5").Value = "https://web.archive.org/web/20200626061624/https://www.moh.gov.bh/COVID19"
$ws.Range("D9").Value = "https://minsalud.gob.bo/4335-covid-19-bolivia-reporta-1-016-nuevos-contagios-y-37-fallecidos"
$ws.Range("D11").Value = "https://web.archive.org/web/20200626061652/https://coronavirus.bg/"
$ws.Range("D12").Value = "https://web.archive.org/web/20200626061653/https://www.canada.ca/en/public-health/services/diseases/2019-novel-coronavirus-infection.html"
$ws.Range("D20").Value = "https://www.gestionderiesgos.gob.ec/wp-content/uploads/2020/06/INFOGRAFIA-NACIONALCOVI-19-COE-NACIONAL-25062020-08h00.pdf"
$ws.Range("D42").Value = "https://www.mhlw.go.jp/stf/newpage_12088.html"
$ws.Range("F42").Value = "See Table: 国内の発生状況, column 1 '検査実施人数'."
$ws.Range("D43").Value = "https://www.mhlw.go.jp/content/10906000/000643532.pdf"
$ws.Range("D46").Value = "https://twitter.com/KUWAIT_MOH/status/1276481026156572673/photo/2"
$ws.Range("D51").Value = "https://www.facebook.com/561317043971945/videos/594103654855023/"
$ws.Range("F51").Value = "Numbers visible in video at time: 1:24"
$ws.Range("D60").Value = "https://twitter.com/OmanVSCovid19/status/1276089447549997056"
$ws.Range("D61").Value = "https://web.archive.org/web/20200626062509/http://www.covid.gov.pk/"
$ws.Range("D64").Value = "https://www.gob.pe/institucion/minsa/noticias/188959-minsa-casos-confirmados-por-coronavirus-covid-19-ascienden-a-268-602-en-el-peru-comunicado-n-146"
$ws.Range("D66").Value = "https://twitter.com/MZ_GOV_PL/status/1276077673983467521"
$ws.Range("D83").Value = "https://www.folkhalsomyndigheten.se/smittskydd-beredskap/utbrott/aktuella-utbrott/covid-19/antal-individer-som-har-testats-for-covid-19/"
$ws.Range("D90").Value = "https://twitter.com/MinofHealthUG/status/1276068542790668288"
$ws.Range("D94").Value = "https://web.archive.org/web/20200626063118/https://www.cdc.gov/coronavirus/2019-ncov/cases-updates/testing-in-us.html"

# --- Numeric updates (dates stored as serials, plus stats columns) ---
$ws.Range("C5").Value = 44008
$ws.Range("G5").Value = 93
$ws.Range("H5").Value = 511458
$ws.Range("I5").Value = 300.578
$ws.Range("J5").Value = 8695
$ws.Range("K5").Value = 5.11
$ws.Range("L5").Value = 7010
$ws.Range("M5").Value = 4.12

$ws.Range("C9").Value = 44007
$ws.Range("G9").Value = 103
$ws.Range("H9").Value = 66597
$ws.Range("I9").Value = 5.705
$ws.Range("J9").Value = 1651
$ws.Range("K9").Value = 0.141
$ws.Range("L9").Value = 1977
$ws.Range("M9").Value = 0.169

$ws.Range("C11").Value = 44008
$ws.Range("G11").Value = 57
$ws.Range("H11").Value = 128293
$ws.Range("I11").Value = 18.464
$ws.Range("J11").Value = 2775
$ws.Range("K11").Value = 0.399
$ws.Range("L11").Value = 1992
$ws.Range("M11").Value = 0.287

$ws.Range("C12").Value = 44008
$ws.Range("G12").Value = 99
$ws.Range("H12").Value = 2558287
$ws.Range("I12").Value = 67.783
$ws.Range("J12").Value = 39713
$ws.Range("K12").Value = 1.052
$ws.Range("L12").Value = 37560
$ws.Range("M12").Value = 0.995

$ws.Range("C14").Value = 44007
$ws.Range("G14").Value = 113
$ws.Range("H14").Value = 670093
$ws.Range("I14").Value = 13.169
$ws.Range("J14").Value = 18501
$ws.Range("K14").Value = 0.364
$ws.Range("L14").Value = 17069
$ws.Range("M14").Value = 0.335

$ws.Range("C15").Value = 44007
$ws.Range("G15").Value = 107
$ws.Range("H15").Value = 28943
$ws.Range("I15").Value = 5.682
$ws.Range("J15").Value = 703
$ws.Range("K15").Value = 0.138
$ws.Range("L15").Value = 498
$ws.Range("M15").Value = 0.098

$ws.Range("C20").Value = 44007
$ws.Range("G20").Value = 96
$ws.Range("H20").Value = 106476
$ws.Range("I20").Value = 6.035
$ws.Range("L20").Value = 1540
$ws.Range("M20").Value = 0.087

$ws.Range("C29").Value = 44005
$ws.Range("G29").Value = 47
$ws.Range("H29").Value = 283124
$ws.Range("I29").Value = 9.112
$ws.Range("J29").Value = 3008
$ws.Range("K29").Value = 0.097
$ws.Range("L29").Value = 3115
$ws.Range("M29").Value = 0.1

$ws.Range("C35").Value = 44008
$ws.Range("G35").Value = 94
$ws.Range("H35").Value = 7776228
$ws.Range("I35").Value = 5.635
$ws.Range("J35").Value = 215446
$ws.Range("K35").Value = 0.156
$ws.Range("L35").Value = 192800
$ws.Range("M35").Value = 0.14

$ws.Range("C42").Value = 44007
$ws.Range("G42").Value = 133
$ws.Range("H42").Value = 435495
$ws.Range("I42").Value = 3.443
$ws.Range("L42").Value = 5580
$ws.Range("M42").Value = 0.044

$ws.Range("C43").Value = 44005
$ws.Range("G43").Value = 39
$ws.Range("H43").Value = 632744
$ws.Range("I43").Value = 5.003
$ws.Range("J43").Value = 5261
$ws.Range("K43").Value = 0.042
$ws.Range("L43").Value = 5780
$ws.Range("M43").Value = 0.046

$ws.Range("C46").Value = 44008
$ws.Range("G46").Value = 45
$ws.Range("H46").Value = 372284
$ws.Range("I46").Value = 87.174
$ws.Range("J46").Value = 3774
$ws.Range("K46").Value = 0.884
$ws.Range("L46").Value = 3267
$ws.Range("M46").Value = 0.765

$ws.Range("C50").Value = 44007
$ws.Range("G50").Value = 107
$ws.Range("H50").Value = 716178
$ws.Range("I50").Value = 22.127
$ws.Range("J50").Value = 11842
$ws.Range("K50").Value = 0.366
$ws.Range("L50").Value = 6656
$ws.Range("M50").Value = 0.206

$ws.Range("C51").Value = 44006
$ws.Range("G51").Value = 84
$ws.Range("H51").Value = 45185
$ws.Range("I51").Value = 83.592
$ws.Range("J51").Value = 1187
$ws.Range("K51").Value = 2.196
$ws.Range("L51").Value = 1099
$ws.Range("M51").Value = 2.033

$ws.Range("C52").Value = 44003
$ws.Range("G52").Value = 173
$ws.Range("H52").Value = 456138
$ws.Range("I52").Value = 3.538
$ws.Range("J52").Value = 2241
$ws.Range("K52").Value = 0.017
$ws.Range("L52").Value = 8215
$ws.Range("M52").Value = 0.064

$ws.Range("C57").Value = 44007
$ws.Range("G57").Value = 109
$ws.Range("H57").Value = 378257
$ws.Range("I57").Value = 78.44
$ws.Range("J57").Value = 9825
$ws.Range("K57").Value = 2.037
$ws.Range("L57").Value = 7257
$ws.Range("M57").Value = 1.505

$ws.Range("C60").Value = 44007
$ws.Range("G60").Value = 22
$ws.Range("J60").Value = 3835
$ws.Range("K60").Value = 0.751
$ws.Range("L60").Value = 3496
$ws.Range("M60").Value = 0.685

$ws.Range("C61").Value = 44008
$ws.Range("G61").Value = 107
$ws.Range("H61").Value = 1193017
$ws.Range("I61").Value = 5.401
$ws.Range("J61").Value = 21041
$ws.Range("K61").Value = 0.095
$ws.Range("L61").Value = 25987
$ws.Range("M61").Value = 0.118

$ws.Range("C64").Value = 44008
$ws.Range("G64").Value = 95
$ws.Range("H64").Value = 232730
$ws.Range("I64").Value = 7.058
$ws.Range("L64").Value = 3360
$ws.Range("M64").Value = 0.102

$ws.Range("C65").Value = 44005
$ws.Range("G65").Value = 82
$ws.Range("H65").Value = 580560
$ws.Range("I65").Value = 5.298
$ws.Range("J65").Value = 11927
$ws.Range("K65").Value = 0.109
$ws.Range("L65").Value = 12197
$ws.Range("M65").Value = 0.111

$ws.Range("C66").Value = 44007
$ws.Range("G66").Value = 59
$ws.Range("H66").Value = 1278454
$ws.Range("I66").Value = 33.78
$ws.Range("J66").Value = 21453
$ws.Range("K66").Value = 0.567
$ws.Range("L66").Value = 18748
$ws.Range("M66").Value = 0.495

$ws.Range("C67").Value = 44007
$ws.Range("L67").Value = 12155
$ws.Range("M67").Value = 0.321

$ws.Range("C69").Value = 44007
$ws.Range("G69").Value = 98
$ws.Range("H69").Value = 337496
$ws.Range("I69").Value = 117.143
$ws.Range("J69").Value = 4324
$ws.Range("K69").Value = 1.501
$ws.Range("L69").Value = 3975
$ws.Range("M69").Value = 1.38

$ws.Range("C73").Value = 44007
$ws.Range("G73").Value = 49
$ws.Range("H73").Value = 1417771
$ws.Range("I73").Value = 40.724
$ws.Range("J73").Value = 37740
$ws.Range("K73").Value = 1.084
$ws.Range("L73").Value = 31357
$ws.Range("M73").Value = 0.901

$ws.Range("C80").Value = 44007
$ws.Range("G80").Value = 117
$ws.Range("H80").Value = 1460012
$ws.Range("I80").Value = 24.617
$ws.Range("J80").Value = 43118
$ws.Range("K80").Value = 0.727
$ws.Range("L80").Value = 33131
$ws.Range("M80").Value = 0.559

$ws.Range("C83").Value = 44003
$ws.Range("G83").Value = 18
$ws.Range("H83").Value = 447462
$ws.Range("I83").Value = 44.306
$ws.Range("L83").Value = 8829
$ws.Range("M83").Value = 0.874

$ws.Range("C90").Value = 44006
$ws.Range("G90").Value = 77
$ws.Range("H90").Value = 153528
$ws.Range("I90").Value = 3.356
$ws.Range("J90").Value = 3455
$ws.Range("K90").Value = 0.076
$ws.Range("L90").Value = 2560
$ws.Range("M90").Value = 0.056

$ws.Range("C94").Value = 44007
$ws.Range("G94").Value = 30
$ws.Range("H94").Value = 31281178
$ws.Range("I94").Value = 94.504
$ws.Range("J94").Value = 1171117
$ws.Range("K94").Value = 3.538
$ws.Range("L94").Value = 682954
$ws.Range("M94").Value = 2.063

$ws.Range("C95").Value = 44007
$ws.Range("G95").Value = 111
$ws.Range("H95").Value = 29207820
$ws.Range("I95").Value = 88.24
$ws.Range("J95").Value = 637587
$ws.Range("K95").Value = 1.926
$ws.Range("L95").Value = 538019
$ws.Range("M95").Value = 1.625

